$wb = $excel.ActiveWorkbook

# --- Add the hidden "DropdownOptions" sheet right after Sheet1 ---
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "DropdownOptions"
$newSheet.Move($null, $wb.Worksheets.Item("Sheet1"))

$ws1 = $wb.Worksheets.Item("Sheet1")
$ddSheet = $wb.Worksheets.Item("DropdownOptions")
$ddSheet.Range("A1").Value = "0% - 10%"
$ddSheet.Range("A2").Value = "11% - 25%"
$ddSheet.Range("A3").Value = "26% - 50%"
$ddSheet.Range("A4").Value = "51% - 75%"
$ddSheet.Range("A5").Value = "76% - 90%"
$ddSheet.Range("A6").Value = "91% - 99%"
# Force text format so "100%" isn't auto-converted into the number 1
$ddSheet.Range("A7").NumberFormat = "@"
$ddSheet.Range("A7").Value = "100%"
$ddSheet.Visible = $false

# --- Add the new "Status as of July 4, 2025" column on Sheet1 ---
$ws1.Range("AH1").Value = "Status as of July 4, 2025"

# Tidy up the stray empty cell at AE2 left over from the old sheet
$ws1.Range("AE2").ClearContents()

# --- Attach a dropdown-list data validation to AH2 sourced from DropdownOptions ---
$validation = $ws1.Range("AH2").Validation
$validation.Add(3, 1, 1, "=DropdownOptions!`$A`$1:`$A`$7")
$validation.IgnoreBlank = $true
$validation.InCellDropdown = $true
$validation.ShowInput = $false
$validation.ShowError = $false
